# Auto-generated Excel COM-interop script implementing the profit-table refresh
# described in the commit diff for Sheets/Excalibur_Profits.xlsx (ALC/ARM/BSM/CRP/CUL/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1490.25
$ws.Range("I15").Value = 1490.25
$ws.Range("K15").Value = 4470.75
$ws.Range("M15").Value = -4301.75
# Row 18
$ws.Range("H18").Value = 6472
$ws.Range("I18").Value = 7962.6665
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 7962.6665
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -7678.6665
$ws.Range("N18").Value = -2568
# Row 76
$ws.Range("H76").Value = 3799.1875
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
# Row 79
$ws.Range("H79").Value = 3799.1875
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
# Row 125
$ws.Range("H125").Value = 2965.75
$ws.Range("I125").Value = 5000
$ws.Range("J125").Value = 2780.818
$ws.Range("K125").Value = 45000
$ws.Range("L125").Value = 25027.362
$ws.Range("M125").Value = -42540
$ws.Range("N125").Value = -29947.362
# Row 137
$ws.Range("H137").Value = 2907.074
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2907.074
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 8721.222
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -13821.222
# Row 138
$ws.Range("H138").Value = 3416.2341
$ws.Range("I138").Value = 2636.4614
$ws.Range("J138").Value = 3714.3823
$ws.Range("K138").Value = 7909.3842
$ws.Range("L138").Value = 11143.1469
$ws.Range("M138").Value = -2769.3842
$ws.Range("N138").Value = -21423.1469
# Row 140
$ws.Range("H140").Value = 130848.75
$ws.Range("J140").Value = 130848.75
$ws.Range("L140").Value = 130848.75
$ws.Range("N140").Value = -141208.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6953255
$ws.Range("I32").Value = 8137863.5
$ws.Range("J32").Value = 14834.429
$ws.Range("K32").Value = 8137863.5
$ws.Range("L32").Value = 14834.429
$ws.Range("M32").Value = -8137576.5
$ws.Range("N32").Value = -15408.429
# Row 61
$ws.Range("H61").Value = 5193.265
$ws.Range("I61").Value = 2274.7878
$ws.Range("K61").Value = 2274.7878
$ws.Range("M61").Value = -2062.7878
# Row 74
$ws.Range("H74").Value = 2953.4578
$ws.Range("I74").Value = 2415.8857
$ws.Range("K74").Value = 2415.8857
$ws.Range("M74").Value = -1541.8857
# Row 77
$ws.Range("H77").Value = 2953.4578
$ws.Range("I77").Value = 2415.8857
$ws.Range("K77").Value = 12079.4285
$ws.Range("M77").Value = -7711.428499999998
# Row 80
$ws.Range("H80").Value = 128286.5
$ws.Range("J80").Value = 129943.8
$ws.Range("L80").Value = 129943.8
$ws.Range("N80").Value = -131939.8
# Row 83
$ws.Range("H83").Value = 128286.5
$ws.Range("J83").Value = 129943.8
$ws.Range("L83").Value = 389831.4
$ws.Range("N83").Value = -399815.4
# Row 132
$ws.Range("H132").Value = 3989.6118
$ws.Range("I132").Value = 2576.392
$ws.Range("K132").Value = 7729.175999999999
$ws.Range("M132").Value = -5199.175999999999
# Row 136
$ws.Range("H136").Value = 5193.265
$ws.Range("I136").Value = 2274.7878
$ws.Range("K136").Value = 6824.3634
$ws.Range("M136").Value = -4274.3634

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2614.818
$ws.Range("I99").Value = 2176.3
$ws.Range("K99").Value = 2176.3
$ws.Range("M99").Value = -678.3000000000002
# Row 105
$ws.Range("H105").Value = 3345.2083
$ws.Range("I105").Value = 2911.5293
$ws.Range("J105").Value = 4398.4287
$ws.Range("K105").Value = 2911.5293
$ws.Range("L105").Value = 4398.4287
$ws.Range("M105").Value = -1164.5293
$ws.Range("N105").Value = -7892.4287

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 508
$ws.Range("J22").Value = 443.2857
$ws.Range("L22").Value = 443.2857
$ws.Range("N22").Value = -1143.2857
# Row 28
$ws.Range("H28").Value = 17500
$ws.Range("J28").Value = 17500
$ws.Range("L28").Value = 17500
$ws.Range("N28").Value = -17990
# Row 31
$ws.Range("H31").Value = 7287.413
$ws.Range("I31").Value = 2282.182
$ws.Range("J31").Value = 8860.485000000001
$ws.Range("K31").Value = 2282.182
$ws.Range("L31").Value = 8860.485000000001
$ws.Range("M31").Value = -1987.182
$ws.Range("N31").Value = -9450.485000000001
# Row 34
$ws.Range("H34").Value = 7287.413
$ws.Range("I34").Value = 2282.182
$ws.Range("J34").Value = 8860.485000000001
$ws.Range("K34").Value = 2282.182
$ws.Range("L34").Value = 8860.485000000001
$ws.Range("M34").Value = -2080.182
$ws.Range("N34").Value = -9264.485000000001
# Row 45
$ws.Range("H45").Value = 16500
$ws.Range("J45").Value = 19000
$ws.Range("L45").Value = 19000
$ws.Range("N45").Value = -20186
# Row 59
$ws.Range("H59").Value = 129997
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 129997
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 129997
$ws.Range("M59").Value = ""
$ws.Range("N59").Value = -132287
# Row 68
$ws.Range("H68").Value = 79412.336
$ws.Range("J68").Value = 79412.336
$ws.Range("L68").Value = 79412.336
$ws.Range("N68").Value = -80910.336
# Row 71
$ws.Range("H71").Value = 79412.336
$ws.Range("J71").Value = 79412.336
$ws.Range("L71").Value = 238237.008
$ws.Range("N71").Value = -245725.008
# Row 74
$ws.Range("H74").Value = 66996
$ws.Range("J74").Value = 89997
$ws.Range("L74").Value = 89997
$ws.Range("N74").Value = -91745
# Row 77
$ws.Range("H77").Value = 66996
$ws.Range("J77").Value = 89997
$ws.Range("L77").Value = 269991
$ws.Range("N77").Value = -278727
# Row 132
$ws.Range("H132").Value = 2773.205
$ws.Range("I132").Value = 1553.6562
$ws.Range("K132").Value = 4660.9686
$ws.Range("M132").Value = -2130.9686
# Row 140
$ws.Range("H140").Value = 98221.89
$ws.Range("J140").Value = 98221.89
$ws.Range("L140").Value = 98221.89
$ws.Range("N140").Value = -108581.89

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 2138.1538
$ws.Range("I113").Value = 1833.3334
$ws.Range("J113").Value = 2177.913
$ws.Range("K113").Value = 5500.0002
$ws.Range("L113").Value = 6533.739
$ws.Range("M113").Value = -3330.0002
$ws.Range("N113").Value = -10873.739
# Row 131
$ws.Range("H131").Value = 13814.444
$ws.Range("J131").Value = 20328.916
$ws.Range("L131").Value = 60986.74800000001
$ws.Range("N131").Value = -71066.74800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1542.0385
$ws.Range("I93").Value = 1383.8636
$ws.Range("J93").Value = 2412
$ws.Range("K93").Value = 1383.8636
$ws.Range("L93").Value = 2412
$ws.Range("M93").Value = -135.8635999999999
$ws.Range("N93").Value = -4908
# Row 100
$ws.Range("H100").Value = 8275.75
$ws.Range("J100").Value = 13412.444
$ws.Range("L100").Value = 13412.444
$ws.Range("N100").Value = -14494.444
# Row 139
$ws.Range("H139").Value = 97853.25
$ws.Range("J139").Value = 97853.25
$ws.Range("L139").Value = 97853.25
$ws.Range("N139").Value = -108133.25

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1423.6666
$ws.Range("I113").Value = 955.5
$ws.Range("K113").Value = 2866.5
$ws.Range("M113").Value = -696.5
